{"js": "// Update the \"About this document\" / Acknowledgments export page:\n//  1. Fix a typo: \"From author annotations\" -> \"For author annotations\"\n//  2. Add two new bullet items (\"Credit Title\", \"Credit Authors\") after the\n//     existing \"Resource Title\" bullet, matching its list level/numbering.\n\nconst body = context.document.body;\n\n// 1) Fix \"From author annotations\" -> \"For author annotations\".\nconst fromResults = body.search(\"From author annotations\", { matchCase: true });\nfromResults.load(\"items\");\nawait context.sync();\n\nif (fromResults.items.length > 0) {\n  fromResults.items[0].insertText(\"For author annotations\", Word.InsertLocation.replace);\n}\n\n// 2) Insert \"Credit Title\" and \"Credit Authors\" bullets right after the\n//    \"Resource Title\" paragraph, inheriting its paragraph formatting\n//    (numbering level/list) automatically via insertParagraph.\nconst titleResults = body.search(\"Resource Title\", { matchCase: true });\ntitleResults.load(\"items\");\nawait context.sync();\n\nif (titleResults.items.length > 0) {\n  const resourceTitlePara = titleResults.items[0].paragraphs.getFirst();\n  const creditTitlePara = resourceTitlePara.insertParagraph(\"Credit Title\", Word.InsertLocation.after);\n  const creditAuthorsPara = creditTitlePara.insertParagraph(\"Credit Authors\", Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "# Update the \"About this document\" / Acknowledgments export page:\n#  1. Fix a typo: \"From author annotations\" -> \"For author annotations\"\n#  2. Add two new bullet items (\"Credit Title\", \"Credit Authors\") after the\n#     existing \"Resource Title\" bullet, matching its list level/numbering.\n\n$d = $word.ActiveDocument\n\n# 1) Fix \"From author annotations\" -> \"For author annotations\".\n$find = $d.Content.Find\n$find.Text = \"From author annotations\"\n$find.Replacement.Text = \"For author annotations\"\n$find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# 2) Insert \"Credit Title\" and \"Credit Authors\" bullets right after the\n#    \"Resource Title\" paragraph, inheriting its paragraph formatting\n#    (numbering level/list) by using InsertParagraphAfter on its Range.\n$targetIndex = -1\n$idx = 0\nforeach ($p in $d.Paragraphs) {\n    $idx = $idx + 1\n    $t = $p.Range.Text\n    $t = $t.TrimEnd([char]13, [char]7)\n    if ($t -eq \"Resource Title\") {\n        $targetIndex = $idx\n        break\n    }\n}\n\nif ($targetIndex -ge 0) {\n    $resourceTitlePara = $d.Paragraphs.Item($targetIndex)\n    $resourceTitlePara.Range.InsertParagraphAfter()\n\n    $creditTitlePara = $d.Paragraphs.Item($targetIndex + 1)\n    $creditTitlePara.Range.Text = \"Credit Title\"\n\n    $creditTitlePara.Range.InsertParagraphAfter()\n    $creditAuthorsPara = $d.Paragraphs.Item($targetIndex + 2)\n    $creditAuthorsPara.Range.Text = \"Credit Authors\"\n}\n"}
